$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (old "NB" / SVM row 9); row 8 (NB) will be overwritten below
$ws.Rows.Item(9).Delete()

# Update header row text (B1:G1 existing, H1:L1 new)
$ws.Range("B1").Value = "Algorithm"
$ws.Range("C1").Value = "One Year Base mean"
$ws.Range("D1").Value = "One Year Base std"
$ws.Range("E1").Value = "Two Year Base mean"
$ws.Range("F1").Value = "Two Year Base std"
$ws.Range("G1").Value = "Three Year Base mean"
$ws.Range("H1").Value = "Three Year Base std"
$ws.Range("I1").Value = "Five Year Base mean"
$ws.Range("J1").Value = "Five Year Base std"
$ws.Range("K1").Value = "Ten Year Base mean"
$ws.Range("L1").Value = "Ten Year Base std"

# Apply header style (bold, centered, bordered) to the newly added header cells H1:L1
$ws.Range("B1").Copy()
$ws.Range("H1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update algorithm names and statistics for each row
$ws.Range("B2").Value = "LR"
$ws.Range("C2").Value = 0.9118301951749158
$ws.Range("D2").Value = 0.009018271541097652
$ws.Range("E2").Value = 0.8983306931030371
$ws.Range("F2").Value = 0.005063791703037574
$ws.Range("G2").Value = 0.8906387395441685
$ws.Range("H2").Value = 0.01120332276323281
$ws.Range("I2").Value = 0.8746928166351606
$ws.Range("J2").Value = 0.01840794193848852
$ws.Range("K2").Value = 0.8561271349394864
$ws.Range("L2").Value = 0.01601016145821086

$ws.Range("B3").Value = "LDA"
$ws.Range("C3").Value = 0.9154095495352348
$ws.Range("D3").Value = 0.009837407096394782
$ws.Range("E3").Value = 0.9035580651938157
$ws.Range("F3").Value = 0.009286384141093978
$ws.Range("G3").Value = 0.8895897884952175
$ws.Range("H3").Value = 0.01067034765304147
$ws.Range("I3").Value = 0.8765860399839605
$ws.Range("J3").Value = 0.01721737714708205
$ws.Range("K3").Value = 0.8644553783508652
$ws.Range("L3").Value = 0.01438732554179628

$ws.Range("B4").Value = "KNN"
$ws.Range("C4").Value = 0.9007669288419269
$ws.Range("D4").Value = 0.008008831983432128
$ws.Range("E4").Value = 0.8861927311337091
$ws.Range("F4").Value = 0.005520331869793017
$ws.Range("G4").Value = 0.8878338824048105
$ws.Range("H4").Value = 0.01126996757370534
$ws.Range("I4").Value = 0.8820748840006873
$ws.Range("J4").Value = 0.01292487365709362
$ws.Range("K4").Value = 0.865883384232553
$ws.Range("L4").Value = 0.01159807445734856

$ws.Range("B5").Value = "DTREE"
$ws.Range("C5").Value = 0.867581367018882
$ws.Range("D5").Value = 0.01281728093558956
$ws.Range("E5").Value = 0.8691626779316494
$ws.Range("F5").Value = 0.01275134323166882
$ws.Range("G5").Value = 0.8589623161427016
$ws.Range("H5").Value = 0.0145507900168874
$ws.Range("I5").Value = 0.8642743598556454
$ws.Range("J5").Value = 0.01972834541350371
$ws.Range("K5").Value = 0.8508901707951588
$ws.Range("L5").Value = 0.01907601546982548

$ws.Range("B6").Value = "RTREE"
$ws.Range("C6").Value = 0.8340626572389503
$ws.Range("D6").Value = 0.0194712610433196
$ws.Range("E6").Value = 0.8260022938775047
$ws.Range("F6").Value = 0.01609607248018595
$ws.Range("G6").Value = 0.813818843153344
$ws.Range("H6").Value = 0.01203607874688118
$ws.Range("I6").Value = 0.803332474079166
$ws.Range("J6").Value = 0.01481176654405744
$ws.Range("K6").Value = 0.7733785770840403
$ws.Range("L6").Value = 0.0261286162292894

$ws.Range("B7").Value = "XTREE"
$ws.Range("C7").Value = 0.9054844945843594
$ws.Range("D7").Value = 0.008474319486313052
$ws.Range("E7").Value = 0.8944535291078294
$ws.Range("F7").Value = 0.009048054767116644
$ws.Range("G7").Value = 0.8853869423046306
$ws.Range("H7").Value = 0.01508821670337607
$ws.Range("I7").Value = 0.8811304204617059
$ws.Range("J7").Value = 0.0149336883252537
$ws.Range("K7").Value = 0.8611265693926026
$ws.Range("L7").Value = 0.01496201866110236

$ws.Range("B8").Value = "SVM"
$ws.Range("C8").Value = 0.9092283043351607
$ws.Range("D8").Value = 0.008860528125732445
$ws.Range("E8").Value = 0.9001865194951199
$ws.Range("F8").Value = 0.007404630607133845
$ws.Range("G8").Value = 0.8916870782457472
$ws.Range("H8").Value = 0.01469997657671607
$ws.Range("I8").Value = 0.8877581342727845
$ws.Range("J8").Value = 0.01686160118615514
$ws.Range("K8").Value = 0.8753970139124533
$ws.Range("L8").Value = 0.01533431647821785

